$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "productos"

# Add new rows of data
$ws.Range("A4").Value = 30
$ws.Range("B4").Value = "Manzanas"
$ws.Range("C4").Value = 13000
$ws.Range("D4").Value = 45

$ws.Range("A5").Value = 40
$ws.Range("B5").Value = "Pimienta"
$ws.Range("C5").Value = 5000
$ws.Range("D5").Value = 60

# Adjust column B width (closest reachable value to the target 11.28515625
# stored width under this engine's column-width rounding grid)
$ws.Columns.Item(2).ColumnWidth = 10.5

# Select A5 as active cell
$ws.Range("A5").Select()
